$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (slopy%, area [km^2], total time [sec], cpu total time [sec])
# appended starting at row 55 through row 77.
$data = @(
    @(37.37, 0.36, 20.14945220947266, 19.46875),
    @(50.14, 0.04000000000000001, 3.571915864944458, 2.828125),
    @(55.75, 0.16, 12.53389739990234, 11.015625),
    @(32.28, 0.16, 9.142158031463623, 7.828125),
    @(27.15, 0.16, 9.069642543792725, 7.25),
    @(42.14, 0.36, 21.91850090026855, 20.765625),
    @(42.14, 0.36, 26.33934664726257, 25.609375),
    @(39.78, 0.16, 11.15482640266418, 9.734375),
    @(100, 0.16, 36.07674241065979, 31.859375),
    @(93.69, 0.16, 56.04600739479065, 52.6875),
    @(98.48999999999999, 0.16, 51.7009365558624, 49.421875),
    @(97.64, 1, 185.6703400611877, 184.015625),
    @(39.78, 0.16, 17.71249747276306, 15.78125),
    @(39.78, 0.16, 30.19679999351501, 26.5625),
    @(39.78, 0.16, 22.50559282302856, 22.046875),
    @(39.78, 0.16, 31.25962972640991, 29.96875),
    @(39.78, 0.16, 20.78266620635986, 19.890625),
    @(39.78, 0.16, 30.54809165000916, 28.859375),
    @(39.78, 0.16, 10.35814690589905, 9.109375),
    @(41.49, 0.16, 10.58113431930542, 9),
    @(38.24, 0.36, 21.57353019714355, 19.59375),
    @(28.18, 0.36, 26.47938179969788, 23.5625),
    @(53.66, 0.36, 27.1538360118866, 24.921875)
)

$startRow = 55
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the view so the newly appended rows are in frame, matching the
# author's final scroll position/selection (row 55 at the top, C65 selected).
$ws.Range("C65").Select()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
